# Translate English source strings to Polish, per translation diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "To confirm your registration, we would require you and one guest of your choice to provide us with:" "Aby potwierdzić rejestrację, wymagamy dostarczenia przez Państwa oraz wybranego gościa:"

Replace-Text "A scanned copy of your international passports" "Zeskanowanej kopii Państwa międzynarodowych paszportów"

Replace-Text "Covid-19 vaccination certificates" "Świadectwa szczepień przeciwko Covid-19"

Replace-Text "Your country manager will be in touch to confirm your booking or request any other relevant details. " "Państwa krajowy menedżer skontaktuje się z Państwem, aby potwierdzić rezerwację lub poprosić o inne istotne szczegóły. "

Replace-Text "Our event package offers you and your guest: " "Nasz pakiet eventowy oferuje Państwu i Państwa gościom: "

Replace-Text "Flight tickets " "Bilety lotnicze "

Replace-Text "Travel insurance " "Ubezpieczenie podróżne "

Replace-Text "Airport – Hotel – Airport transfer " "Lotnisko - Hotel - Transfer lotniskowy "

Replace-Text "One hotel room for you and your guest / Two hotel rooms for you and your guest" "Jeden pokój hotelowy dla Państwa i Państwa gościa / Dwa pokoje hotelowe dla Państwa i Państwa gościa"

Replace-Text "Meals (Breakfast, lunch, and dinner)" "Posiłki (śniadanie, obiad i kolacja)"

Replace-Text "We will send you a confirmation letter before your departure date with the event agenda and information about your flights, transportation, and accommodation. " "Przed datą wyjazdu wyślemy Państwu list potwierdzający z programem wydarzenia oraz informacjami na temat przelotów, transportu i zakwaterowania. "

Replace-Text "We look forward to seeing you soon." "Czekamy na Państwa wkrótce."
